$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "bell_state.py"
$ws.Range("B9").Value = "Bell State"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "Bell State entangles 2 qubits "

$ws.Range("A10").Value = "ghz_state.py"
$ws.Range("B10").Value = "GHZ State"
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "GHZ State entangles 3 qubits"
